# Re-process the parameter/description cell (B1) the same way the
# "ExcelToCWMS" import routine does: trim stray whitespace from the cell
# text before looking for the "{units=" marker used to extract the unit
# of measure. Trimming here removes the trailing space that had been left
# on the "02600.Flow..." description, matching the cleaned-up value that
# is now stored after parsing throws when "{units=" cannot be located.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$paramCell = $ws.Range("B1")
$paramText = $paramCell.Value()

if ($null -ne $paramText) {
    $trimmedText = $paramText.Trim()

    if ($trimmedText.IndexOf("{units=") -lt 0) {
        throw "Could not find {units= in '$trimmedText'"
    }

    if ($trimmedText -ne $paramText) {
        $paramCell.Value = $trimmedText
    }
}

# Move the active selection to C10, as left by the editor after the change.
$ws.Range("C10").Select()
